# Commit message: "add the NA's under duplicate_image_filename"
#
# The "duplicate_image_filename" header lives in column E (E1).
# The stimuli data rows run from row 2 through row 21 (the practice
# trials p1-p4 in rows 2-5, plus the 16 generic/unique trials in rows
# 6-21). None of those rows had a value in column E yet, so we fill
# them all in with "NA".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E21").Value = "NA"
